$d = $word.ActiveDocument

# Locate the Subtitle paragraph ("Life-sentenced prisoners' experiences of
# parole decision-making") so we can insert the new Author paragraph right
# after it (and before the Date paragraph), matching the target diff.
$subtitlePara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Style.NameLocal -eq "Subtitle") {
        $subtitlePara = $p
        break
    }
}

if ($subtitlePara -ne $null) {
    # Insert a brand-new paragraph immediately after the subtitle.
    $subtitlePara.Range.InsertParagraphAfter()

    # The freshly-created paragraph inherits the Subtitle style and sits
    # right after $subtitlePara - find it by walking forward from the
    # subtitle paragraph.
    $authorPara = $subtitlePara.Next()

    # Give it the document's existing "Author" paragraph style and the
    # author's name as its text.
    $authorPara.Style = "Author"
    $authorPara.Range.Text = "Ben Jarman"
}
